# Weekly data refresh: insert the newest week's price record for
# Hortaliza, Femacal de La Calera - Zapallo italiano.
# This pushes all existing records (previously rows 233-288) down by one
# row, and inserts the new record at row 233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 233, shifting existing data down.
$ws.Rows.Item(233).Insert()

# Populate the new row with the latest weekly record.
$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value = 44543
$ws.Cells.Item(233, 5).Value = 5
$ws.Cells.Item(233, 6).Value = 100112032
$ws.Cells.Item(233, 7).Value = "Zapallo italiano"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 195
$ws.Cells.Item(233, 11).Value = 7500
$ws.Cells.Item(233, 12).Value = 8000
$ws.Cells.Item(233, 13).Value = 7751
$ws.Cells.Item(233, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(233, 15).Value = "Limache"
$ws.Cells.Item(233, 16).Value = 129
$ws.Cells.Item(233, 17).Value = 60
$ws.Cells.Item(233, 18).Value = "Hortaliza"
